$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ "E"=3; "G"=31.910369; "H"=95.73110699999999; "I"=0.579978174461379; "J"=0.579978174461379; "K"=3; "M"=2.211866666666667; "N"=6.6356; "O"=0.4811217919111272; "P"=0.4811217919111272; "Q"=70.58148151213334; "R"=635.2333336092; "S"=0.279040138566203; "T"=0.279040138566203 }
    3 = @{ "E"=3; "G"=31.910369; "H"=95.73110699999999; "I"=0.579978174461379; "J"=0.579978174461379; "K"=3; "M"=2.385444666666666; "N"=7.156333999999999; "O"=0.5188782080888727; "P"=0.5188782080888727; "Q"=76.12041954241532; "R"=685.0837758817379; "S"=0.300938035895176; "T"=0.300938035895176 }
    4 = @{ "E"=3; "G"=11.420477; "H"=34.261431; "I"=0.2075697527013294; "J"=0.2075697527013294; "K"=3; "M"=2.211866666666667; "N"=6.6356; "O"=0.4811217919111272; "P"=0.4811217919111272; "Q"=25.26057239373334; "R"=227.3451515436; "S"=0.09986633136621312; "T"=0.09986633136621312 }
    5 = @{ "E"=3; "G"=11.420477; "H"=34.261431; "I"=0.2075697527013294; "J"=0.2075697527013294; "K"=3; "M"=2.385444666666666; "N"=7.156333999999999; "O"=0.5188782080888727; "P"=0.5188782080888727; "Q"=27.24291595043933; "R"=245.186243553954; "S"=0.1077034213351162; "T"=0.1077034213351162 }
    6 = @{ "E"=3; "G"=0.6836493333333333; "H"=2.050948; "I"=0.01242548126969028; "J"=0.01242548126969028; "K"=3; "M"=2.211866666666667; "N"=6.6356; "O"=0.4811217919111272; "P"=0.4811217919111272; "Q"=1.512141172088889; "R"=13.6092705488; "S"=0.005978169813831537; "T"=0.005978169813831537 }
    7 = @{ "E"=3; "G"=0.6836493333333333; "H"=2.050948; "I"=0.01242548126969028; "J"=0.01242548126969028; "K"=3; "M"=2.385444666666666; "N"=7.156333999999999; "O"=0.5188782080888727; "P"=0.5188782080888727; "Q"=1.630807656070222; "R"=14.677268904632; "S"=0.006447311455858746; "T"=0.006447311455858746 }
    8 = @{ "E"=3; "G"=0.4491346666666667; "H"=1.347404; "I"=0.00816312415756312; "J"=0.00816312415756312; "K"=3; "M"=2.211866666666667; "N"=6.6356; "O"=0.4811217919111272; "P"=0.4811217919111272; "Q"=0.9934259980444445; "R"=8.940833982400001; "S"=0.003927456922279779; "T"=0.003927456922279779 }
    9 = @{ "E"=3; "G"=0.4491346666666667; "H"=1.347404; "I"=0.00816312415756312; "J"=0.00816312415756312; "K"=3; "M"=2.385444666666666; "N"=7.156333999999999; "O"=0.5188782080888727; "P"=0.5188782080888727; "Q"=1.071385895215111; "R"=9.642473056936; "S"=0.004235667235283341; "T"=0.004235667235283341 }
    10 = @{ "E"=3; "G"=0.6568320000000001; "H"=1.970496; "I"=0.01193807017047708; "J"=0.01193807017047708; "K"=3; "M"=2.211866666666667; "N"=6.6356; "O"=0.4811217919111272; "P"=0.4811217919111272; "Q"=1.4528248064; "R"=13.0754232576; "S"=0.005743665712380708; "T"=0.005743665712380708 }
    11 = @{ "E"=3; "G"=0.6568320000000001; "H"=1.970496; "I"=0.01193807017047708; "J"=0.01193807017047708; "K"=3; "M"=2.385444666666666; "N"=7.156333999999999; "O"=0.5188782080888727; "P"=0.5188782080888727; "Q"=1.566836391296; "R"=14.101527521664; "S"=0.006194404458096371; "T"=0.006194404458096371 }
    12 = @{ "E"=3; "G"=9.899486; "H"=29.698458; "I"=0.1799253972395612; "J"=0.1799253972395612; "K"=3; "M"=2.211866666666667; "N"=6.6356; "O"=0.4811217919111272; "P"=0.4811217919111272; "Q"=21.89634310053333; "R"=197.0670879048; "S"=0.08656602953021905; "T"=0.08656602953021905 }
    13 = @{ "E"=3; "G"=9.899486; "H"=29.698458; "I"=0.1799253972395612; "J"=0.1799253972395612; "K"=3; "M"=2.385444666666666; "N"=7.156333999999999; "O"=0.5188782080888727; "P"=0.5188782080888727; "Q"=23.61467608144133; "R"=212.532084732972; "S"=0.09335936770934213; "T"=0.09335936770934213 }
}

foreach ($rowKey in $updates.Keys) {
    $rowData = $updates[$rowKey]
    foreach ($col in $rowData.Keys) {
        $cellRef = "$col$rowKey"
        $ws.Range($cellRef).Value = $rowData[$col]
    }
}
